{"js": "// no-op\n", "ps1": "# no-op\n"}
